$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5630,45744),
    @(5590,45744.01041666666),
    @(5550,45744.02083333334),
    @(5510,45744.03125),
    @(5490,45744.04166666666),
    @(5460,45744.05208333334),
    @(5440,45744.0625),
    @(5430,45744.07291666666),
    @(5430,45744.08333333334),
    @(5440,45744.09375),
    @(5440,45744.10416666666),
    @(5450,45744.11458333334),
    @(5460,45744.125),
    @(5470,45744.13541666666),
    @(5490,45744.14583333334),
    @(5520,45744.15625),
    @(5560,45744.16666666666),
    @(5620,45744.17708333334),
    @(5700,45744.1875),
    @(5800,45744.19791666666),
    @(5920,45744.20833333334),
    @(6060,45744.21875),
    @(6210,45744.22916666666),
    @(6370,45744.23958333334),
    @(6540,45744.25),
    @(6700,45744.26041666666),
    @(6840,45744.27083333334),
    @(6970,45744.28125),
    @(7080,45744.29166666666),
    @(7170,45744.30208333334),
    @(7220,45744.3125),
    @(7250,45744.32291666666),
    @(7260,45744.33333333334),
    @(7250,45744.34375),
    @(7220,45744.35416666666),
    @(7170,45744.36458333334),
    @(7130,45744.375),
    @(7080,45744.38541666666),
    @(7030,45744.39583333334),
    @(6980,45744.40625),
    @(6940,45744.41666666666),
    @(6910,45744.42708333334),
    @(6880,45744.4375),
    @(6860,45744.44791666666),
    @(6840,45744.45833333334),
    @(6820,45744.46875),
    @(6810,45744.47916666666),
    @(6800,45744.48958333334),
    @(6790,45744.5),
    @(6790,45744.51041666666),
    @(6800,45744.52083333334),
    @(6810,45744.53125),
    @(6820,45744.54166666666),
    @(6840,45744.55208333334),
    @(6850,45744.5625),
    @(6860,45744.57291666666),
    @(6880,45744.58333333334),
    @(6890,45744.59375),
    @(6920,45744.60416666666),
    @(6940,45744.61458333334),
    @(6980,45744.625),
    @(7020,45744.63541666666),
    @(7070,45744.64583333334),
    @(7120,45744.65625),
    @(7170,45744.66666666666),
    @(7230,45744.67708333334),
    @(7270,45744.6875),
    @(7320,45744.69791666666),
    @(7380,45744.70833333334),
    @(7430,45744.71875),
    @(7500,45744.72916666666),
    @(7580,45744.73958333334),
    @(7660,45744.75),
    @(7730,45744.76041666666),
    @(7770,45744.77083333334),
    @(7780,45744.78125),
    @(7760,45744.79166666666),
    @(7720,45744.80208333334),
    @(7640,45744.8125),
    @(7550,45744.82291666666),
    @(7440,45744.83333333334),
    @(7330,45744.84375),
    @(7190,45744.85416666666),
    @(7060,45744.86458333334),
    @(6910,45744.875),
    @(6760,45744.88541666666),
    @(6620,45744.89583333334),
    @(6460,45744.90625),
    @(6310,45744.91666666666),
    @(6170,45744.92708333334),
    @(6030,45744.9375),
    @(5920,45744.94791666666),
    @(5830,45744.95833333334),
    @(5770,45744.96875),
    @(5720,45744.97916666666),
    @(5660,45744.98958333334)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
